$d = $word.ActiveDocument

# The paragraph currently reads "Payment Mode: Cashfree / Cash / Other" and
# needs to become "Payment Mode: {{payment_mode}}", keeping the run / proof
# error layout that Word would produce for a manually-typed edit (new runs
# for newly typed text, the spell-check wrapper kept around
# "payment"/"_"/"mode").

$full = $d.Content.Text
$base = $full.IndexOf("Payment Mode: Cash")

$cashStart = $base + "Payment Mode: ".Length
$cashEnd = $cashStart + "Cash".Length
$freeStart = $cashEnd
$freeEnd = $freeStart + "free".Length
$spaceStart = $freeEnd
$spaceEnd = $spaceStart + " ".Length
$tailStart = $spaceEnd
$tailEnd = $tailStart + "/ Cash / Other".Length

# --- Phase 1: make every text-content change. Each freshly written bit of
# text is marked Bold=true right away; this stops the engine from silently
# re-merging it into an identically-formatted neighbouring run while we are
# still editing. The Bold marker is cleared afterwards, in Phase 2.

# 1) "Cash" -> "payment"
$rCash = $d.Range($cashStart, $cashEnd)
$rCash.Font.Bold = $true
$rCash.Text = "payment"
$paymentStart = $cashStart
$paymentEnd = $cashStart + "payment".Length

$shift1 = "payment".Length - "Cash".Length
$freeStart += $shift1; $freeEnd += $shift1
$spaceStart += $shift1; $spaceEnd += $shift1
$tailStart += $shift1; $tailEnd += $shift1

# 2) "free" -> "_"
$rFree = $d.Range($freeStart, $freeEnd)
$rFree.Font.Bold = $true
$rFree.Text = "_"
$underStart = $freeStart
$underEnd = $freeStart + "_".Length

$shift2 = "_".Length - "free".Length
$spaceStart += $shift2; $spaceEnd += $shift2
$tailStart += $shift2; $tailEnd += $shift2

# 3) insert "mode" right after "_" (still inside the spellStart/spellEnd
#    proofErr pair)
$rModeIns = $d.Range($underEnd, $underEnd)
$rModeIns.InsertBefore("mode")
$modeStart = $underEnd
$modeEnd = $modeStart + "mode".Length
$d.Range($modeStart, $modeEnd).Font.Bold = $true

$shift3 = "mode".Length
$spaceStart += $shift3; $spaceEnd += $shift3
$tailStart += $shift3; $tailEnd += $shift3

# 4) the old " " run becomes the closing "}}"
$rSpace = $d.Range($spaceStart, $spaceEnd)
$rSpace.Font.Bold = $true
$rSpace.Text = "}}"
$closeStart = $spaceStart
$closeEnd = $spaceStart + "}}".Length

$shift4 = "}}".Length - " ".Length
$tailStart += $shift4; $tailEnd += $shift4

# 5) drop the old "/ Cash / Other" run entirely
$d.Range($tailStart, $tailEnd).Text = ""

# 6) insert the opening "{{" right before the spellStart-wrapped "payment"
$rOpenIns = $d.Range($paymentStart, $paymentStart)
$rOpenIns.InsertBefore("{{")
$openStart = $paymentStart
$openEnd = $paymentStart + 2
$d.Range($openStart, $openEnd).Font.Bold = $true

# Everything from $paymentStart onward shifted right by 2 because "{{" was
# inserted in front of it.
$paymentStart += 2; $paymentEnd += 2
$underStart += 2; $underEnd += 2
$modeStart += 2; $modeEnd += 2
$closeStart += 2; $closeEnd += 2

# --- Phase 2: clear the Bold marker, one run's exact range at a time. These
# are pure formatting edits (no text content change), so they do not trigger
# another run-coalescing pass, and the separate runs set up above survive.
$d.Range($openStart, $openEnd).Font.Bold = $false
$d.Range($paymentStart, $paymentEnd).Font.Bold = $false
$d.Range($underStart, $underEnd).Font.Bold = $false
$d.Range($modeStart, $modeEnd).Font.Bold = $false
$d.Range($closeStart, $closeEnd).Font.Bold = $false
